# Parts list update: the previously-placeholder "Voltage regulator" row (row 3)
# has been identified as a diode (D4), replacing it with the real part info.
# Row 2 (the e-stop MOSFET) is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: was blank / "Voltage regulator 5V@0.5A" / "AZ2940D-5.0TRE1DICT-ND"
# now identified as the diode D4.
$ws.Range("A3").Value = "D4"
$ws.Range("B3").Value = "Diode"
$ws.Range("C3").Value = "VSKY20301608-G4-08GICT-ND"

# Widen column C slightly so the longer Digikey part numbers fit.
$ws.Columns.Item(3).ColumnWidth = 24.78

# Move the active selection down to B5.
$ws.Range("B5").Select()
